# Commit: "novas traducoes e textos" (new translations and texts)
#
# This adds a new localization row for the "sword_" key with version
# "EA 23.117" right before the existing "SpMeteor" row (i.e. the new row
# becomes row 15 on the "Calc" sheet), pushing every following row down
# by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; this shifts rows 15.. down to 16..
$ws.Rows("15").Insert()

# Populate the newly inserted row with the new id/version pair.
$ws.Range("A15").Value = "sword_"
$ws.Range("B15").Value = "EA 23.117"

# Restore the selection over the data range, matching the saved view.
[void]$ws.Range("A3:B35").Select()
